$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$meta = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty, now "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was a "Contact" row ("Contact" / "No display for ContactDetail");
# it becomes the "Jurisdiction" row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" row ("Contact" / "No display for ContactDetail")
# and is removed entirely, shifting all following rows up by one.
$meta.Rows.Item(11).Delete()

# --- Sheet 2: "Elements" ---
$elements = $wb.Worksheets.Item(2)

# Row 2 (the root Extension element) gets a real Short/Definition instead of the
# generic "Extension" / "An Extension" placeholders.
$elements.Range("K2").Value = "Copay Exclusion"
$elements.Range("L2").Value = "Customer-specific code for the exclusion reason of a copayment"
